# Add the new "ProductosBusqueda" sheet after the last existing sheet
# (so it lands at the end of the tab strip, becomes sheetId=3 / rId3,
# and is activated - matching activeTab moving to the 3rd tab).
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "ProductosBusqueda"

# Column A entries first (this is the order the original author typed
# them in, which drives the shared-string table ordering).
$ws.Range("A2").Value = "Phones"
$ws.Range("A3").Value = "Tablets"
$ws.Range("A4").Value = "Cameras"

# Header row.
$ws.Range("A1").Value = "Categoria"
$ws.Range("B1").Value = "SubCategoria"
$ws.Range("C1").Value = "Producto"
$ws.Range("D1").Value = "Cantidad"

# Product names.
$ws.Range("C2").Value = "iPhone"
$ws.Range("C3").Value = "Samsung Galaxy Tab 10.1"
$ws.Range("C4").Value = "Canon EOS 5D"

# Quantities to add to cart.
$ws.Range("D2").Value = 2
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1

# A trailing styled-but-empty helper cell (underline/black "Aptos Narrow"
# font - a new, distinct style from the existing hyperlink style).
$ws.Range("E5").Font.Underline = $true

# Leave the selection parked at F12, as in the authored workbook.
$ws.Range("F12").Select() | Out-Null
